$wb = $excel.ActiveWorkbook

# --- Sheet 1: data_src_a_table_a (Table2) ---
$ws1 = $wb.Worksheets.Item(1)
$lo1 = $ws1.ListObjects.Item(1)
$lo1Name = $lo1.Name
$lo1.Unlist()
$ws1.Columns.Item(7).Insert()
$ws1.Range("G3").Value = "created_at_field"
$ws1.Range("H3").Value = "updated_at_field"
$newLo1 = $ws1.ListObjects.Add(1, $ws1.Range("B3:N7"), $null, 1, $null)
$newLo1.Name = $lo1Name

# --- Sheet 2: data_src_a_table_b (Table24) ---
$ws2 = $wb.Worksheets.Item(2)
$lo2 = $ws2.ListObjects.Item(1)
$lo2Name = $lo2.Name
$lo2.Unlist()
$ws2.Columns.Item(7).Insert()
$ws2.Range("G3").Value = "created_at_field"
$ws2.Range("H3").Value = "updated_at_field"
$newLo2 = $ws2.ListObjects.Add(1, $ws2.Range("B3:N6"), $null, 1, $null)
$newLo2.Name = $lo2Name

# --- Sheet view / tab selection changes ---
# sheet2 becomes the active / selected sheet
$ws1.Select()
$ws1Sel = $ws1.Range("D27")
$ws2.Select()
$ws2Sel = $ws2.Range("E13")
